$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new team's players (columns C, D, E, F)
$ws.Range("C2").Value = "Andrea Conzatti | FC Savignano"
$ws.Range("D2").Value = "Leonardo  Parisi  | MediaserT"
$ws.Range("E2").Value = "marco bertolini | Fc Wanda Tim"
$ws.Range("F2").Value = "Matteo  Tatarella | Bayern Mona"

# Remove rows 3 and 4 entirely (old data no longer needed)
$ws.Rows("3:4").Delete()
